$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Row 3: the old "Position av det intressanta ljudet" label in D3 is removed ---
$ws.Range("D3").ClearContents()

# --- Row 4 ---
$ws.Range("E4").Value = 23
$ws.Range("F4").Value = 3
$ws.Range("H4").Value = 1

# --- Row 5 ---
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 1

# --- Row 6 ---
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 10
$ws.Range("I6").Value = 2

# --- Row 7 ---
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 13
$ws.Range("I7").Value = 4

# --- Row 8 ---
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 18

# --- Build a clustered-column chart summarising "Slagljud utan filter" ---
$co = $ws.ChartObjects().Add(500, 30, 430, 330)
$co.Name = "Diagram 3"
$chart = $co.Chart
$chart.ChartType = 51

$cols = @("E", "F", "G", "H", "I", "J")
$names = @("Placering 1 på testbänk", "Placering 2 på testbänk", "Placering 3 på testbänk", "Placering 4 på testbänk", "Placering 5 på testbänk", "Ej placerbar")

for ($i = 0; $i -lt 6; $i++) {
    $col = $cols[$i]
    $series = $chart.SeriesCollection().NewSeries()
    $series.Name = $names[$i]
    $series.XValues = "=Blad1!`$D`$4:`$D`$8"
    $series.Values = "=Blad1!`$" + $col + "`$4:`$" + $col + "`$8"
    $series.HasDataLabels = $true
    $labels = $series.DataLabels()
    $labels.ShowValue = $true
    $labels.ShowCategoryName = $false
    $labels.ShowSeriesName = $false
    $labels.ShowPercentage = $false
    $labels.ShowLegendKey = $false
    $labels.ShowBubbleSize = $false
    $labels.Position = -4165
}

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Slagljud utan filter"

$chart.Axes(1).HasTitle = $true
$chart.Axes(1).AxisTitle.Text = "Var slagljudet placerar på testbänken"

$chart.Axes(2).HasTitle = $true
$chart.Axes(2).AxisTitle.Text = "Antal val från testarna"

$chart.HasLegend = $true
$chart.Legend.Position = -4152

# --- Restore the active selection to D4, as in the edited workbook ---
$ws.Range("D4").Select()
